$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for D (Price) and E (Volume) columns so that values
# like "1.00" or "600.60" are preserved as literal text instead of being
# coerced into numbers by Excel.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '67.975.26'
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").Value = '3.681.55'
$ws.Range("E3").Value = '  -0.14%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '600.60'
$ws.Range("E5").Value = '  +4.09%  '
$ws.Range("D6").Value = '190.23'
$ws.Range("E6").Value = '  +11.08%  '
$ws.Range("D7").Value = '0.625'
$ws.Range("E7").Value = '  +0.42%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = '0.710'
$ws.Range("E9").Value = '  +1.34%  '
$ws.Range("D10").Value = '58.18'
$ws.Range("E10").Value = '  +14.26%  '
$ws.Range("E11").Value = '  -3.37%  '
$ws.Range("E12").Value = '  -3.21%  '
$ws.Range("D13").Value = '10.24'
$ws.Range("E13").Value = '  -0.86%  '
$ws.Range("D14").Value = '4.270.84'
$ws.Range("E14").Value = '  -0.04%  '
$ws.Range("D15").Value = '3.686.20'
$ws.Range("E15").Value = '  +0.24%  '
$ws.Range("E16").Value = '  +0.89%  '
$ws.Range("D17").Value = '19.06'
$ws.Range("E17").Value = '  -1.14%  '
$ws.Range("E18").Value = '  +1.70%  '
$ws.Range("D19").Value = '12.60'
$ws.Range("E19").Value = '  -1.15%  '
$ws.Range("D20").Value = '67.908.14'
$ws.Range("E20").Value = '  +0.71%  '
$ws.Range("D21").Value = '402.82'
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("D22").Value = '4.47'
$ws.Range("E22").Value = '  +0.19%  '
$ws.Range("D23").Value = '88.57'
$ws.Range("E23").Value = '  +1.15%  '
$ws.Range("D24").Value = '11.55'
$ws.Range("E24").Value = '  +7.57%  '
$ws.Range("E25").Value = '  -0.92%  '
$ws.Range("D26").Value = '12.66'
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("D27").Value = '6.04'
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("D28").Value = '3.71'
$ws.Range("E28").Value = '  -1.62%  '
$ws.Range("D29").Value = '9.39'
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("D30").Value = '32.07'
$ws.Range("E30").Value = '  -1.04%  '
$ws.Range("D31").Value = '7.64'
$ws.Range("E31").Value = '  +3.65%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").Value = '45.78'
$ws.Range("E32").Value = '  +6.92%  '
$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D33").Value = '12.42'
$ws.Range("E33").Value = '  +0.73%  '
$ws.Range("E34").Value = '  +3.85%  '
$ws.Range("E35").Value = '  +2.87%  '
$ws.Range("D36").Value = '620.21'
$ws.Range("E36").Value = '  +2.39%  '
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("E38").Value = '  +2.57%  '
$ws.Range("D39").Value = '0.0₃0788'
$ws.Range("E39").Value = '  -10.38%  '
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("D41").Value = '0.135'
$ws.Range("E41").Value = '  +1.04%  '
$ws.Range("E42").Value = '  -1.21%  '
$ws.Range("E43").Value = '  -0.15%  '
$ws.Range("D44").Value = '2.58'
$ws.Range("E44").Value = '  -6.28%  '
$ws.Range("D45").Value = '2.852.04'
$ws.Range("E45").Value = '  +1.29%  '
$ws.Range("E46").Value = '  +3.28%  '
$ws.Range("D47").Value = '3.24'
$ws.Range("E47").Value = '  +3.78%  '
$ws.Range("D48").Value = '9.04'
$ws.Range("E48").Value = '  -1.35%  '
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").Value = '144.95'
$ws.Range("E49").Value = '  +4.43%  '
$ws.Range("B50").Value = 'WEMIXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D50").Value = '2.65'
$ws.Range("E50").Value = '  -1.13%  '
$ws.Range("D51").Value = '2.55'
$ws.Range("E51").Value = '  -10.61%  '

# Restore the original (default/General) style so no stray style index
# is left attached to the edited cells.
$dataRange.Style = "Normal"

